$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

function Set-ParaText($paraIndex, $newText) {
    $para = $tr.Paragraphs($paraIndex, 1)
    $chars = $para.Characters(1, $para.Length)
    $chars.Text = $newText
}

# Paragraph 1: "View client list" -> "Delete clients, orders, vendor inventory"
Set-ParaText 1 "Delete clients, orders, vendor inventory"

# Paragraph 2: "View vendor info" -> "View client info, totals, subtotals, items not found"
Set-ParaText 2 "View client info, totals, subtotals, items not found"

# Paragraph 3: "View expense and items not found report" -> "Read in excel clients to database"
Set-ParaText 3 "Read in excel clients to database"

# Paragraph 4: "Read in excel clients to database" -> "Read in test excel vendor to database"
Set-ParaText 4 "Read in test excel vendor to database"

# Paragraph 5: "Read in excel vendors to database" -> two runs: "Error " + "handling for adding clients"
Set-ParaText 5 "Error handling for adding clients"
$para5 = $tr.Paragraphs(5, 1)
$firstPart = $para5.Characters(1, 6)
$firstPart.Text = "Error "

# Delete the old paragraph 6 ("Delete all records from all tables") and the old
# paragraph 7 ("Error handling for adding clients") entirely -- their content was
# folded into paragraph 5 above, so both whole paragraphs go away. Re-fetch
# paragraph index 6 fresh each time since the collection shifts after a delete.
$toDelete = $tr.Paragraphs(6, 1)
$toDelete.Delete()

$toDelete2 = $tr.Paragraphs(6, 1)
$toDelete2.Delete()
